$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-11 first so the remaining rows 2-4 can be
# updated in place to become the new (smaller) data set.
$ws.Rows("5:11").Delete()

# Ensure date-like text values stay as plain text (not auto-converted
# to Excel date serials) for column A.
$ws.Range("A2:A4").NumberFormat = "@"

# Row 2: 2025-03-08, Vasanth Kumar, 09:22:02, (blank exit time)
$ws.Range("A2").Value = "2025-03-08"
$ws.Range("B2").Value = "Vasanth Kumar"
$ws.Range("C2").Value = "09:22:02"
$ws.Range("D2").ClearContents()

# Row 3: 2025-03-08, Sreenath (unchanged), 09:57:16, nan
$ws.Range("A3").Value = "2025-03-08"
$ws.Range("C3").Value = "09:57:16"
$ws.Range("D3").Value = "nan"

# Row 4: 2025-03-08, Soundharraja, 10:00:24, nan
$ws.Range("A4").Value = "2025-03-08"
$ws.Range("B4").Value = "Soundharraja"
$ws.Range("C4").Value = "10:00:24"
$ws.Range("D4").Value = "nan"
